# Update gh-pages output (杭州-漫展信息.xlsx)
# Applies numeric "want-to-go" / "min price" corrections to a handful of
# existing events, and inserts a brand-new event row
# ("杭州·百鬼夜行·咒术回战only") into the two sheets that list it
# ("展览" and "全部类型"), shifting the following rows down by one.

$wb = $excel.ActiveWorkbook

function Update-EventValues {
    param($ws, $rowMap)
    foreach ($row in $rowMap.Keys) {
        $vals = $rowMap[$row]
        if ($vals.ContainsKey('F')) {
            $ws.Cells.Item($row, 6).Value = $vals['F']
        }
        if ($vals.ContainsKey('G')) {
            $ws.Cells.Item($row, 7).Value = $vals['G']
        }
    }
}

function Insert-NewEvent {
    param($ws, $insertRow, $lastRowBefore)

    # Push rows [$insertRow .. end] down by one row.
    $ws.Rows.Item($insertRow).Insert()

    # Copy the index column's number format/border style down into the
    # freshly created (blank) row so A$insertRow matches its neighbours.
    $ws.Cells.Item($insertRow - 1, 1).Copy()
    $ws.Cells.Item($insertRow, 1).PasteSpecial(-4122)

    # Fill in the new event's data.
    $ws.Cells.Item($insertRow, 2).Value = "2024-03-16"
    $ws.Cells.Item($insertRow, 3).Value = "杭州·百鬼夜行·咒术回战only"
    $ws.Cells.Item($insertRow, 4).Value = "长生路18号 梅地亚宾馆"
    $ws.Cells.Item($insertRow, 5).Value = "2024.03.16 09:00-03.16 17:00"
    $ws.Cells.Item($insertRow, 6).Value = 0
    $ws.Cells.Item($insertRow, 7).Value = 79
    $ws.Cells.Item($insertRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81478"
    $ws.Cells.Item($insertRow, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/4weHdCdk1706495040356.jpeg"

    # The "index" column (A) is always just (row number - 1). Re-number
    # it for the new row and for every row that just shifted down.
    for ($r = $insertRow; $r -le ($lastRowBefore + 1); $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) -- rows 1-based, header on row 1
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

Update-EventValues $wsExpo @{
    3  = @{ F = 7663; G = 80 }
    8  = @{ F = 585 }
    13 = @{ F = 3071 }
    16 = @{ F = 721 }
    19 = @{ F = 450 }
    21 = @{ F = 216 }
    22 = @{ F = 215 }
    23 = @{ F = 274 }
    24 = @{ F = 284 }
    27 = @{ F = 255 }
}

Insert-NewEvent $wsExpo 28 35

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) -- rows 1-based, header on row 1
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

Update-EventValues $wsAll @{
    6  = @{ F = 7663; G = 80 }
    11 = @{ F = 585 }
    17 = @{ F = 3071 }
    21 = @{ F = 721 }
    25 = @{ F = 450 }
    27 = @{ F = 216 }
    28 = @{ F = 215 }
    29 = @{ F = 274 }
    30 = @{ F = 284 }
    33 = @{ F = 255 }
}

Insert-NewEvent $wsAll 34 42
